$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 16

    # Trade #
    $ws.Cells.Item($row, 1).Value = 15

    # Date / Time - must stay plain text, not be auto-converted to a date/time
    # serial by Excel's smart input parsing, so force the cell to Text first.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = "20:03:34"
    $ws.Cells.Item($row, 3).Style = "Normal"

    # Strategy
    $ws.Cells.Item($row, 4).Value = "MarketMaking"

    # Side
    $ws.Cells.Item($row, 5).Value = "UP"

    # Entry Price
    $ws.Cells.Item($row, 6).Value = 0.86

    # Exit Price - trade is still OPEN, so no exit price yet
    $ws.Cells.Item($row, 7).Value = ""

    # Status
    $ws.Cells.Item($row, 8).Value = "OPEN"

    # P&L %
    $ws.Cells.Item($row, 9).Value = 0

    # P&L $
    $ws.Cells.Item($row, 10).Value = 0

    # Capital After
    $ws.Cells.Item($row, 11).Value = 99.68557117791565

    # Entry Slippage (bps)
    $ws.Cells.Item($row, 12).Value = 0

    # Exit Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0

    # Confidence
    $ws.Cells.Item($row, 14).Value = 0.6

    # Entry Reason
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"

    # Exit Reason - trade is still OPEN, so no exit reason yet
    $ws.Cells.Item($row, 16).Value = ""

    # Duration (min)
    $ws.Cells.Item($row, 17).Value = 0
}
